$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the galaxy name in A8: "J1341-0321" -> "J1341+0321"
$ws.Range("A8").Value = "J1341+0321"

# Update the active selection to reflect where work left off
$ws.Range("A8").Select()
